$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").Value = "25.589.30"
$ws.Range("E2").Value = "  -5.96%  "

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").Value = "1.803.11"
$ws.Range("E3").Value = "  -5.25%  "

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5: 'BNB' -> 'BNB'
$ws.Range("D5").Value = "272.99"
$ws.Range("E5").Value = "  -10.77%  "

# Row 6: 'USDC' -> 'USDC'
$ws.Range("E6").Value = "  +0.06%  "

# Row 7: 'XRP' -> 'XRP'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5010"
$ws.Range("E7").Value = "  -7.27%  "

# Row 8: 'Cardano' -> 'Cardano'
$ws.Range("D8").Value = "0.3495"
$ws.Range("E8").Value = "  -8.12%  "

# Row 9: 'OKB' -> 'Dogecoin'
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.06571"
$ws.Range("E9").Value = "  -9.88%  "

# Row 10: 'Dogecoin' -> 'Solana'
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "19.84"
$ws.Range("E10").Value = "  -10.00%  "

# Row 11: 'Solana' -> 'Polygon'
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "0.8318"
$ws.Range("E11").Value = "  -7.79%  "

# Row 12: 'Polygon' -> 'TRON'
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07764"
$ws.Range("E12").Value = "  -5.12%  "

# Row 13: 'TRON' -> 'WrappedEther'
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.807.59"
$ws.Range("E13").Value = "  +39.00%  "

# Row 14: 'WrappedEther' -> 'Polkadot'
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.036"
$ws.Range("E14").Value = "  -5.79%  "

# Row 15: 'Polkadot' -> 'Litecoin'
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "87.17"
$ws.Range("E15").Value = "  -8.73%  "

# Row 16: 'Litecoin' -> 'BinanceUSD'
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.04%  "

# Row 17: 'BinanceUSD' -> 'Avalanche'
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "13.86"
$ws.Range("E17").Value = "  -6.40%  "

# Row 18: 'Avalanche' -> 'Dai'
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.08%  "

# Row 19: 'Dai' -> 'ShibaInu'
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007923"
$ws.Range("E19").Value = "  -8.19%  "

# Row 20: 'ShibaInu' -> 'WrappedBTC'
$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").Value = "25.657.77"
$ws.Range("E20").Value = "  -5.83%  "

# Row 21: 'WrappedBTC' -> 'Uniswap'
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "4.699"
$ws.Range("E21").Value = "  -6.86%  "

# Row 22: 'Uniswap' -> 'Cosmos'
$ws.Range("B22").Value = "Cosmos"
$ws.Range("C22").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D22").Value = "9.968"
$ws.Range("E22").Value = "  -7.77%  "

# Row 23: 'Cosmos' -> 'Chainlink'
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "6.033"
$ws.Range("E23").Value = "  -7.31%  "

# Row 24: 'Chainlink' -> 'Monero'
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").Value = "141.74"
$ws.Range("E24").Value = "  -4.39%  "

# Row 25: 'Monero' -> 'LidoDAOToken'
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").Value = "2.101"
$ws.Range("E25").Value = "  -8.79%  "

# Row 26: 'LidoDAOToken' -> 'Toncoin'
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "1.651"
$ws.Range("E26").Value = "  -5.91%  "

# Row 27: 'Toncoin' -> 'EthereumClassic'
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "16.87"
$ws.Range("E27").Value = "  -8.09%  "

# Row 28: 'EthereumClassic' -> 'BitcoinCash'
$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D28").Value = "107.92"
$ws.Range("E28").Value = "  -7.45%  "

# Row 29: 'BitcoinCash' -> 'InternetComputer(DFINITY)'
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "4.301"
$ws.Range("E29").Value = "  -11.39%  "

# Row 30: 'InternetComputer(DFINITY)' -> 'Filecoin'
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "4.176"
$ws.Range("E30").Value = "  -10.23%  "

# Row 31: 'Filecoin' -> 'Stellar'
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.08754"
$ws.Range("E31").Value = "  -4.87%  "

# Row 32: 'Stellar' -> 'Hedera'
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.04768"
$ws.Range("E32").Value = "  -5.74%  "

# Row 33: 'Hedera' -> 'HuobiToken'
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "2.878"
$ws.Range("E33").Value = "  -4.42%  "

# Row 34: 'HuobiToken' -> 'ARBITRUM'
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "1.127"
$ws.Range("E34").Value = "  -7.63%  "

# Row 35: 'ARBITRUM' -> 'ImmutableX'
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.7169"
$ws.Range("E35").Value = "  -12.93%  "

# Row 36: 'ImmutableX' -> 'Frax'
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("E36").Value = "  -0.26%  "

# Row 37: 'Frax' -> 'MXToken'
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "3.015"
$ws.Range("E37").Value = "  -9.03%  "

# Row 38: 'MXToken' -> 'VeChain'
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01853"
$ws.Range("E38").Value = "  -7.31%  "

# Row 39: 'VeChain' -> 'TheSandbox'
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "0.5115"
$ws.Range("E39").Value = "  -14.85%  "

# Row 40: 'TheSandbox' -> 'RenderToken'
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "2.248"
$ws.Range("E40").Value = "  -16.74%  "

# Row 41: 'RenderToken' -> 'TrustWalletToken'
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "0.9422"
$ws.Range("E41").Value = "  -12.26%  "

# Row 42: 'TrustWalletToken' -> 'Quant'
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "113.38"
$ws.Range("E42").Value = "  -2.14%  "

# Row 43: 'Quant' -> 'FraxShare'
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "6.125"
$ws.Range("E43").Value = "  -7.90%  "

# Row 44: 'FraxShare' -> 'Aptos'
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "7.949"
$ws.Range("E44").Value = "  -14.08%  "

# Row 45: 'Aptos' -> 'PaxDollar'
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.09%  "

# Row 46: 'PaxDollar' -> 'Algorand'
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "0.1371"
$ws.Range("E46").Value = "  -10.33%  "

# Row 47: 'Algorand' -> 'Decentraland'
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.4526"
$ws.Range("E47").Value = "  -12.26%  "

# Row 48: 'Decentraland' -> 'EnergySwap'
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.194"
$ws.Range("E48").Value = "  -9.75%  "

# Row 49: 'EnergySwap' -> 'Elrond'
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "35.88"
$ws.Range("E49").Value = "  -5.76%  "

# Row 50: 'Elrond' -> 'NEARProtocol'
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "1.478"
$ws.Range("E50").Value = "  -9.71%  "

# Row 51: 'NEARProtocol' -> 'Cronos'
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.05776"
$ws.Range("E51").Value = "  -5.31%  "
